$d = $word.ActiveDocument

# Colors (Word COM uses BGR-packed integers, i.e. 0x<BB><GG><RR>)
$BLACK    = 0x000000   # target color "000000"
$GREEN    = 0x33A900   # target color "00A933" (RGB 00,A9,33 -> BGR 33,A9,00)

# ---------------------------------------------------------------------------
# Helper: recolor an entire paragraph (text runs + paragraph mark) to $color
# ---------------------------------------------------------------------------
function Set-ParaColor($paraIndex, $color) {
    $p = $d.Paragraphs.Item($paraIndex)
    $p.Range.Font.Color = $color
}

# 1) "What is Hypervisor?" Q&A  -> color 00A933 -> 000000
Set-ParaColor 3 $BLACK
Set-ParaColor 4 $BLACK

# 2) "What is Virtualization?" Q&A -> unchanged (still 00A933)

# 3) "What is Containerization?" Q&A -> color 00A933 -> 000000
Set-ParaColor 7 $BLACK

# Merge the split runs of the Containerization answer before recoloring
# (the three trailing runs hold identical formatting, so replacing the
# concatenated text with itself collapses them into a single run).
$containerizationAnswer = " Containerization bundles applications, dependencies, and configurations into containers. It solves compatibility issues when deploying software on different machines."
$rng8 = $d.Paragraphs.Item(8).Range
$rng8.Find.Execute($containerizationAnswer, $true, $false, $false, $false, $false, $true, 1, $false, $containerizationAnswer, 2)
Set-ParaColor 8 $BLACK

# 4) "Difference between Virtualization and Containerization:" Q&A -> 00A933 -> 000000
Set-ParaColor 9 $BLACK
Set-ParaColor 10 $BLACK

# 5) "What is Docker?" Q&A -> 00A933 -> 000000
Set-ParaColor 11 $BLACK

# Merge the split runs of the Docker answer before recoloring
$dockerAnswer = " Docker is a containerization platform that packages applications and dependencies into containers. It solves compatibility issues when deploying software on different machines."
$rng12 = $d.Paragraphs.Item(12).Range
$rng12.Find.Execute($dockerAnswer, $true, $false, $false, $false, $false, $true, 1, $false, $dockerAnswer, 2)
Set-ParaColor 12 $BLACK

# 6) "What are Docker Images?" Q&A -> 00A933 -> 000000
Set-ParaColor 13 $BLACK
Set-ParaColor 14 $BLACK

# 7) "What is Docker Hub?" Q&A -> 00A933 -> 000000
Set-ParaColor 15 $BLACK
Set-ParaColor 16 $BLACK

# 8) "Explain Docker Architecture:" -> "EXplain Docker Architecture:" + text color 00A933
#    (paragraph-mark stays uncolored, so only the text itself is recolored).
#    The first two letters are recolored as a separate operation from the
#    remainder, mirroring how the original edit (capitalizing "Ex" -> "EX")
#    split the run in two.
$rng17 = $d.Paragraphs.Item(17).Range
$rng17.Find.Execute("Explain Docker Architecture:", $true, $false, $false, $false, $false, $true, 1, $false, "EXplain Docker Architecture:", 2)

$p17 = $d.Paragraphs.Item(17)
$p17Start = $p17.Range.Start
$rngEX = $d.Range($p17Start, $p17Start + 2)
$rngEX.Font.Color = $GREEN
$rngRest = $d.Range($p17Start + 2, $p17.Range.End - 1)
$rngRest.Font.Color = $GREEN

# 9) Docker Architecture answer paragraph -> color none -> 00A933 (incl. paragraph mark)
Set-ParaColor 18 $GREEN

# 10) "What is a Dockerfile?" Q&A -> 00A933 -> 000000
Set-ParaColor 19 $BLACK
Set-ParaColor 20 $BLACK

# 11) "Tell us something about Docker Compose:" Q&A -> 00A933 -> 000000
Set-ParaColor 21 $BLACK
Set-ParaColor 22 $BLACK
